# Auto-generated data refresh for Unintended_Deviation worksheet
# Shifts data window from 20-21 Aug 2025 to 17-18 Sep 2025, refreshes values,
# and appends 2 additional quarter-hour rows (134, 135).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, 45917.04166666666, 0, 16.201, 1, "17.09.20251"),
  @(3, 45917.05208333334, 0, 9.157, 2, "17.09.20252"),
  @(4, 45917.0625, 0, 5.422, 3, "17.09.20253"),
  @(5, 45917.07291666666, 0, 0.242, 4, "17.09.20254"),
  @(6, 45917.08333333334, 0, 18.496, 5, "17.09.20255"),
  @(7, 45917.09375, 2.023, 0, 6, "17.09.20256"),
  @(8, 45917.10416666666, 13.027, 0, 7, "17.09.20257"),
  @(9, 45917.11458333334, 8.142, 0, 8, "17.09.20258"),
  @(10, 45917.125, 5.948, 0, 9, "17.09.20259"),
  @(11, 45917.13541666666, 2.622, 0, 10, "17.09.202510"),
  @(12, 45917.14583333334, 5.048, 0, 11, "17.09.202511"),
  @(13, 45917.15625, 28.006, 0, 12, "17.09.202512"),
  @(14, 45917.16666666666, 8.198, 0, 13, "17.09.202513"),
  @(15, 45917.17708333334, 16.857, 0, 14, "17.09.202514"),
  @(16, 45917.1875, 10.706, 0, 15, "17.09.202515"),
  @(17, 45917.19791666666, 4.864, 0, 16, "17.09.202516"),
  @(18, 45917.20833333334, 2.534, 0, 17, "17.09.202517"),
  @(19, 45917.21875, 8.988, 0, 18, "17.09.202518"),
  @(20, 45917.22916666666, 7.457, 0, 19, "17.09.202519"),
  @(21, 45917.23958333334, 7.323, 0, 20, "17.09.202520"),
  @(22, 45917.25, 9.591, 0, 21, "17.09.202521"),
  @(23, 45917.26041666666, 17.734, 0, 22, "17.09.202522"),
  @(24, 45917.27083333334, 2.661, 0, 23, "17.09.202523"),
  @(25, 45917.28125, 10.051, 0, 24, "17.09.202524"),
  @(26, 45917.29166666666, 22.571, 0, 25, "17.09.202525"),
  @(27, 45917.30208333334, 41.02, 0, 26, "17.09.202526"),
  @(28, 45917.3125, 24.796, 0, 27, "17.09.202527"),
  @(29, 45917.32291666666, 0, 3.261, 28, "17.09.202528"),
  @(30, 45917.33333333334, 8.995, 0, 29, "17.09.202529"),
  @(31, 45917.34375, 23.091, 0, 30, "17.09.202530"),
  @(32, 45917.35416666666, 14.633, 0, 31, "17.09.202531"),
  @(33, 45917.36458333334, 4.184, 0, 32, "17.09.202532"),
  @(34, 45917.375, 0, 9.768, 33, "17.09.202533"),
  @(35, 45917.38541666666, 0, 1.938, 34, "17.09.202534"),
  @(36, 45917.39583333334, 0, 7.104, 35, "17.09.202535"),
  @(37, 45917.40625, 0, 12.038, 36, "17.09.202536"),
  @(38, 45917.41666666666, 0, 3.262, 37, "17.09.202537"),
  @(39, 45917.42708333334, 0, 8.342, 38, "17.09.202538"),
  @(40, 45917.4375, 3.516, 0, 39, "17.09.202539"),
  @(41, 45917.44791666666, 3.067, 0, 40, "17.09.202540"),
  @(42, 45917.45833333334, 24.97, 0, 41, "17.09.202541"),
  @(43, 45917.46875, 5.466, 0, 42, "17.09.202542"),
  @(44, 45917.47916666666, 0.77, 0, 43, "17.09.202543"),
  @(45, 45917.48958333334, 0, 9.234, 44, "17.09.202544"),
  @(46, 45917.5, 4.943, 0, 45, "17.09.202545"),
  @(47, 45917.51041666666, 0, 66.209, 46, "17.09.202546"),
  @(48, 45917.52083333334, 0, 28.481, 47, "17.09.202547"),
  @(49, 45917.53125, 0, 28.767, 48, "17.09.202548"),
  @(50, 45917.54166666666, 0, 21.878, 49, "17.09.202549"),
  @(51, 45917.55208333334, 0, 21.549, 50, "17.09.202550"),
  @(52, 45917.5625, 0, 55.442, 51, "17.09.202551"),
  @(53, 45917.57291666666, 0, 52.114, 52, "17.09.202552"),
  @(54, 45917.58333333334, 2.416, 0, 53, "17.09.202553"),
  @(55, 45917.59375, 0.911, 0, 54, "17.09.202554"),
  @(56, 45917.60416666666, 10.987, 0, 55, "17.09.202555"),
  @(57, 45917.61458333334, 0, 0.415, 56, "17.09.202556"),
  @(58, 45917.625, 0, 51.64, 57, "17.09.202557"),
  @(59, 45917.63541666666, 0, 30.924, 58, "17.09.202558"),
  @(60, 45917.64583333334, 0, 4.305, 59, "17.09.202559"),
  @(61, 45917.65625, 4.188, 0, 60, "17.09.202560"),
  @(62, 45917.66666666666, 0, 19.371, 61, "17.09.202561"),
  @(63, 45917.67708333334, 0, 12.639, 62, "17.09.202562"),
  @(64, 45917.6875, 0, 19.967, 63, "17.09.202563"),
  @(65, 45917.69791666666, 0, 20.564, 64, "17.09.202564"),
  @(66, 45917.70833333334, 0, 9.183, 65, "17.09.202565"),
  @(67, 45917.71875, 1.529, 0, 66, "17.09.202566"),
  @(68, 45917.72916666666, 0, 15.215, 67, "17.09.202567"),
  @(69, 45917.73958333334, 0, 30.322, 68, "17.09.202568"),
  @(70, 45917.75, 0, 13.383, 69, "17.09.202569"),
  @(71, 45917.76041666666, 0, 20.818, 70, "17.09.202570"),
  @(72, 45917.77083333334, 0, 30.109, 71, "17.09.202571"),
  @(73, 45917.78125, 0, 5.728, 72, "17.09.202572"),
  @(74, 45917.79166666666, 14.003, 0, 73, "17.09.202573"),
  @(75, 45917.80208333334, 16.67, 0, 74, "17.09.202574"),
  @(76, 45917.8125, 9.397, 0, 75, "17.09.202575"),
  @(77, 45917.82291666666, 23.502, 0, 76, "17.09.202576"),
  @(78, 45917.83333333334, 53.613, 0, 77, "17.09.202577"),
  @(79, 45917.84375, 23.352, 0, 78, "17.09.202578"),
  @(80, 45917.85416666666, 0, 6.843, 79, "17.09.202579"),
  @(81, 45917.86458333334, 0, 1.048, 80, "17.09.202580"),
  @(82, 45917.875, 0, 8.29, 81, "17.09.202581"),
  @(83, 45917.88541666666, 0, 16.524, 82, "17.09.202582"),
  @(84, 45917.89583333334, 0, 12.337, 83, "17.09.202583"),
  @(85, 45917.90625, 0, 20.224, 84, "17.09.202584"),
  @(86, 45917.91666666666, 0, 4.245, 85, "17.09.202585"),
  @(87, 45917.92708333334, 0, 4.371, 86, "17.09.202586"),
  @(88, 45917.9375, 0, 0.909, 87, "17.09.202587"),
  @(89, 45917.94791666666, 0, 1.164, 88, "17.09.202588"),
  @(90, 45917.95833333334, 17.395, 0, 89, "17.09.202589"),
  @(91, 45917.96875, 8.351, 0, 90, "17.09.202590"),
  @(92, 45917.97916666666, 0, 4.895, 91, "17.09.202591"),
  @(93, 45917.98958333334, 0, 5.562, 92, "17.09.202592"),
  @(94, 45918, 0, 4.068, 93, "18.09.202593"),
  @(95, 45918.01041666666, 0, 20.508, 94, "18.09.202594"),
  @(96, 45918.02083333334, 0, 12.698, 95, "18.09.202595"),
  @(97, 45918.03125, 0, 15.691, 96, "18.09.202596"),
  @(98, 45918.04166666666, 0, 18.554, 1, "18.09.20251"),
  @(99, 45918.05208333334, 0, 4.921, 2, "18.09.20252"),
  @(100, 45918.0625, 0, 5.036, 3, "18.09.20253"),
  @(101, 45918.07291666666, 0, 6.034, 4, "18.09.20254"),
  @(102, 45918.08333333334, 0, 10.001, 5, "18.09.20255"),
  @(103, 45918.09375, 4.659, 0, 6, "18.09.20256"),
  @(104, 45918.10416666666, 0, 1.491, 7, "18.09.20257"),
  @(105, 45918.11458333334, 1.387, 0, 8, "18.09.20258"),
  @(106, 45918.125, 0, 0.224, 9, "18.09.20259"),
  @(107, 45918.13541666666, 0, 5.399, 10, "18.09.202510"),
  @(108, 45918.14583333334, 4.209, 0, 11, "18.09.202511"),
  @(109, 45918.15625, 0, 6.711, 12, "18.09.202512"),
  @(110, 45918.16666666666, 0, 9.347, 13, "18.09.202513"),
  @(111, 45918.17708333334, 0, 7.267, 14, "18.09.202514"),
  @(112, 45918.1875, 0, 2.538, 15, "18.09.202515"),
  @(113, 45918.19791666666, 0, 4.577, 16, "18.09.202516"),
  @(114, 45918.20833333334, 0, 31.461, 17, "18.09.202517"),
  @(115, 45918.21875, 0, 14.68, 18, "18.09.202518"),
  @(116, 45918.22916666666, 2.489, 0, 19, "18.09.202519"),
  @(117, 45918.23958333334, 0, 3.721, 20, "18.09.202520"),
  @(118, 45918.25, 0, 2.16, 21, "18.09.202521"),
  @(119, 45918.26041666666, 5.858, 0, 22, "18.09.202522"),
  @(120, 45918.27083333334, 0, 1.462, 23, "18.09.202523"),
  @(121, 45918.28125, 1.094, 0, 24, "18.09.202524"),
  @(122, 45918.29166666666, 7.218, 0, 25, "18.09.202525"),
  @(123, 45918.30208333334, 0, 13.446, 26, "18.09.202526"),
  @(124, 45918.3125, 0, 10.867, 27, "18.09.202527"),
  @(125, 45918.32291666666, 0, 4.542, 28, "18.09.202528"),
  @(126, 45918.33333333334, 0, 1.289, 29, "18.09.202529"),
  @(127, 45918.34375, 0, 17.762, 30, "18.09.202530"),
  @(128, 45918.35416666666, 0, 28.454, 31, "18.09.202531"),
  @(129, 45918.36458333334, 0, 50.526, 32, "18.09.202532"),
  @(130, 45918.375, 0, 30.959, 33, "18.09.202533"),
  @(131, 45918.38541666666, 0, 25.903, 34, "18.09.202534"),
  @(132, 45918.39583333334, 0, 6.328, 35, "18.09.202535"),
  @(133, 45918.40625, 0, 1.228, 36, "18.09.202536"),
  @(134, 45918.41666666666, 4.805, 0, 37, "18.09.202537"),
  @(135, 45918.42708333334, 0, 7.977, 38, "18.09.202538")
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}